# "Menu para administrador y para usuario"
# Adds a "CABECERA TICKETS" / "ESTADOS TICKETS" field-reference block
# (columns I & K, rows 1-10) and a "DETALLE TICKETS" field-reference
# block (column I, rows 19-25) to the existing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Header($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $r.Value = $text
    $r.Font.Bold = $true
}

function Set-Plain($rangeAddr, $text) {
    $ws.Range($rangeAddr).Value = $text
}

# --- CABECERA TICKETS (col I, rows 1-8) ------------------------------------
Set-Header "I1" "CABECERA TICKETS"
Set-Plain "I2" "Id"
Set-Plain "I3" "Fecha Creación"
Set-Plain "I4" "Usuario Creación"
Set-Plain "I5" "Empresa"
Set-Plain "I6" "Título"
Set-Plain "I7" "Descripción"
Set-Plain "I8" "Estado"

# --- ESTADOS TICKETS (col K) ----------------------------------------------
Set-Header "K1" "ESTADOS TICKETS"
Set-Plain "K2" "Borrador"
Set-Plain "K3" "Enviado"
Set-Plain "K4" "Rechazado"
Set-Plain "K5" "Devuelto Incompleto"
Set-Plain "K6" "Resuelto"

# --- CABECERA TICKETS (col I, rows 9-10) ------------------------------------
Set-Plain "I9" "Fecha Estado"
Set-Plain "I10" "Usuario Estado"

# --- DETALLE TICKETS (col I) ---------------------------------------------
Set-Header "I19" "DETALLE TICKETS"
Set-Plain "I20" "Id"
Set-Plain "I21" "IdCabecera"
Set-Plain "I22" "Imagen"
Set-Plain "I23" "Comentario"
Set-Plain "I24" "Usuario"
Set-Plain "I25" "Fecha"

# --- page setup -----------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- selection --------------------------------------------------------
$ws.Range("I26").Select() | Out-Null
